$d = $word.ActiveDocument

# Paragraph 1 (index 0): update date in title line
$old0 = '⚡️🚀המאמר היומי של מייק 14.08.24: ⚡️🚀'
$new0 = '⚡️🚀המאמר היומי של מייק 13.08.24: ⚡️🚀'
$d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2) | Out-Null

# Paragraph 2 (index 1): new paper title + manual line break (same run)
$old1 = 'Jumping Ahead: Improving Reconstruction Fidelity with JumpReLU Sparse Autoencoders'
$new1 = 'Gemma Scope: Open Sparse Autoencoders Everywhere All At Once on Gemma 2' + "^l"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Paragraph 3 (index 2): replace paragraph text
$old2 = 'אתמול סקרנו מאמר שהשתמש בגישת SAE או Sparse AutoEncoders כדי לחדור ל״מחשובותיו״ של מודל שפה גדול דרך האקטיבציות של הנוירונים שלהם. הנחת היסוד במאמר היתה כי נוירונים ״מגיבים״ לכמה קונספטים שונים וניתן לאמן SAE רדוד מאוד (שכבה אחת בדקודר ושכבה אחת באנדוקר) כדי להגיע לוקטור דליל המקודד (נדלק) קונספט אחד בלבד כלומר disentanglement של הפיצ''רים לנוירונים ייעודיים.'
$new2 = 'בזמן האחרון התחלתי להתעניין בשיטות interpretability של מודלי שפה גדולים בעקבות כמה בלוגים מאוד מעניינים של אנטרופיק, OpenAI ולאחר מכן גוגל בנושא הזה. המטרה כאן היא לשפוך קצת אור על הקופסא השחורה שנקראת LLM - הרי אנחנו לא באמת מבינים איך הם עובדים ומה גורם להם לפלוט תשובה כזו אור אחרת לפרומפט שלנו.'
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# Paragraph 4 (index 3): replace paragraph text
$old3 = 'כמאמר יש באנקודר של SAE שכבה לינארית אחת עם פונקציית אקטיבציה הנקראת JumpReLU שראיתי אותה בפעם הראשונה במאמר הזה. פונקציה הזו היא בעצם הזזה של ReLU בציר X ובציר y בפרמטר t נלמד (במאמר זה נקרא טטה). הטענה במאמר שזה מאפשר ללמוד את הייצוג הדליל של דאטה על ידי האנקודר יותר טוב של פונקציית ReLU בגלל שהוא מאפשר לאפס את הקטיבציות בצורה ״נלמדת יותר מ-ReLU".'
$new3 = 'אז המאמר הזה חוקר אחת השיטות המנסות להבין איך מודל שפה מייצג קונספטים סמנטיים שונים. המאמר עושה זאת דרך חקר של אקטיבציות הנוירונים בשכבותיהם השונות של מודלי שפה. עקב כך שיטה זו משויכת למשפחת שיטות המכונות mechanistic interpretability. הרעיון שהמאמר דן בו נקרא SAE או Sparse AutoEncoders. '
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# Paragraph 5 (index 4): replace paragraph text
$old4 = 'עכשיו נשאלת השאלה איך אנחנו אוכפים דלילות על ייצוג הדאטה (אחרי האנקודר). בעבודות קודמות השתמשו ב-L1 בשביל כך אך כאן המחברים משתמשים באותה JumpReLU כדי להפוך את איפוס האיברים בייצוג יותר נלמד. ושימו לב ש- JumpReLU בא עם פרמטר נלמד הזה לזה של האנקודר עצמו שזה עוזר לאכוף דלילות על הייצוג.'
$new4 = 'אז מה הרעיון העיקרי ב- SAE? אנו מנסים להציג אקטיבציות של שכבה מסוימת של LLM על יד וקטור ארוך הרבה יותר מווקטור האקטיבציות אך מאוד דליל. כלומר וקטור n-ממדי של האקטיבציות אנו מייצגים (עם SAE) עם וקטור באורך M >> n אך בווקטור האורך הזה יש פחות מ- n איברים לא שווים לאפס (דלילות). SAE במקרה הזה פשוט מאוד: שכבה אחת לינארית עם אקטיבציה לא לינארית באנקודר (של SAE) ושכבה אחת של דקודר. המטרה כמובן לאמן את SAE כך שיהיה ניתן לשחזר את האקטיבציות המקוריות מייצוגם (אחרי האנקודר).'
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# Paragraph 6 (index 5): replace paragraph text
$old5 = 'יש עוד טריק אחד קטן ולא מאוד מהותי במאמר הנקרא Kernel density estimation או KDE. אם אתם זוכרים KDE עוזר לנו לשערך(כלומר לקרב) פונקצית צפיפות בהינתם דאטהסט של נקודות באמצעות פונקציית קרנל. פונקציית קרנל יכולה להיות גאוסית למשל ומטרתה לשערך את פונקציית הצפיפות לנקודות לא ידועות על ידי קירובה בין הנקודות בדאטהסט (בדומה לספליין). אז המחברים משתמשים בטריק הזה כדי לשערך את JumpReLU בנקודה t שבה היא לא גזירה. '
$new5 = 'אבל למה זה בכלל חשוב ואיך זה קשור ל-interpretability של LLMs. הנחת מוצא של גישה זו (הבלוג של אנטרופיק מדבר על זה בהרחבה) שכל נוירון (או קבוצת נוירונים) בשכבה (מסוימת) הוא ״נדלק״ (מקבל ערכים) על כמה קונספטים לא קשורים (נגיד כלב, מכונה וערפל). כלומר הוא סוג של תערובת עבור כמה קונספטים. אז הייצוג המופק על ידי SAE הוא למעשה מהווה ייצוג של כל קונספט (disentangled). כלומר עבור כל קונספט המקודד קבוצות נוירונים שונות בוקטור הדליל הזה.'
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# Paragraph 7 (index 6): replace paragraph text
$old6 = 'מאמר נחמד בנושא די חשוב שאמשיך לסקור כנראה גם בעתיד…'
$new6 = 'אז מה המאמר הזה עושה? הוא מנסה לאתר שכבות שבהם SAE מאומן עם שגיאת שחזור מינימלית (עם רגולריזציה מתאימה) כלומר הוא מנסה להבין איזו שכבה ב-LLM (וגם בשכבות הפנימיות של בלוקי הטרנספורמר) מקודדת הכי טוב את הקונספטים הסמנטיים.'
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# Paragraph 8 (index 7): replace arXiv pdf link line with lead-in sentence
$old7 = 'https://arxiv.org/pdf/2407.14435'
$new7 = 'בימים הקרובים עוד כמה סקירות בנושא המרתק הזה.'
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# New paragraph (index 8): append new arXiv abs link as its own paragraph
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLastPara.Range.Text = 'https://arxiv.org/abs/2408.05147'

Write-Host "Final paragraph count:" $d.Paragraphs.Count

